$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New StatQuery text shared by both the CasesTab (row 2) and the new
# FilesTab (row 3) rows: combined Trials/Cases/Files counter query. ---
$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
   WHERE c.disease =  "Adenocarcinoma of the colon"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# --- Updated CasesTab main query (row 2, column B) ---
$casesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
 WHERE c.disease = "Adenocarcinoma of the colon"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

# --- New FilesTab main query (row 3, column B) ---
$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
 WHERE c.disease = "Adenocarcinoma of the colon"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

$neo4jFile = $ws.Range("D2").Text
$webFile = $ws.Range("E2").Text

# --- Row 2 (CasesTab): refresh the StatQuery column with the new combined
# Trials/Cases/Files counter query ---
$ws.Range("C2").Value = $statQuery

# --- Row 3 (new FilesTab row) ---
$ws.Range("A3").Value = "FilesTab"

# --- Row 2 (CasesTab): replace the main query ---
$ws.Range("B2").Value = $casesQuery

# --- Row 3 (FilesTab): main query, shared StatQuery, file names ---
$ws.Range("B3").Value = $filesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

[void]$ws.Range("C2").Select()
